# Applies the "checked Tables for BCNF" edit:
#  - Adds a new "Fakultaet" attribute row to the Studiengang block on the
#    "Attributes" sheet (row 27), shifting subsequent rows down.
#  - Updates selections / active-tab bookkeeping to match what Excel
#    records after this interactive edit.
#  - Row-height autofit changes on the "Entity" sheet caused by the edit.

$wb = $excel.ActiveWorkbook

$wsEntity = $wb.Worksheets.Item("Entity")
$wsRelations = $wb.Worksheets.Item("Relations")
$wsAttributes = $wb.Worksheets.Item("Attributes")

# --- Attributes sheet: insert a new row for "Fakultaet" under "Studiengang" ---
$insertRow = $wsAttributes.Rows.Item(27)
$insertRow.Insert()

# A28 (old A27, same entity block marker column) already carries the
# correct "entity marker" style (s="3") - copy it onto the new A27.
$wsAttributes.Cells.Item(28, 1).Copy()
$wsAttributes.Cells.Item(27, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsAttributes.Cells.Item(27, 2).Value = "Fakultät"
$wsAttributes.Cells.Item(27, 3).Value = "Zum Studiengang zugeordnete Fakultät"
$wsAttributes.Cells.Item(27, 4).Value = "String (230)"
$wsAttributes.Cells.Item(27, 5).Value = "not nullable"
$wsAttributes.Cells.Item(27, 6).Value = "no"
$wsAttributes.Cells.Item(27, 7).Value = "no"
$wsAttributes.Cells.Item(27, 8).Value = "no"

# --- Entity sheet: row heights for some description rows return to default ---
$wsEntity.Rows.Item(2).AutoFit()
$wsEntity.Rows.Item(5).AutoFit()
$wsEntity.Rows.Item(6).AutoFit()
$wsEntity.Rows.Item(10).AutoFit()

# --- Selection / active sheet bookkeeping ---
$wsRelations.Range("E10").Select() | Out-Null
$wsAttributes.Range("I27").Select() | Out-Null

$wsAttributes.Activate()
